# Updated symbol list on Tue Jan 17 16:55:03 UTC 2023 with GitHub Actions
#
# Refresh the Price (D) / Volume(1h) (E) columns on the active sheet with the
# latest coinranking.com snapshot. Values are kept as plain text (same as the
# original cells) by formatting the target range as Text ("@") before writing
# the new value, so numeric-looking strings such as "302.59" or "1.77%" are
# not silently re-interpreted as numbers/percentages by Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "302.59" }
    @{ Cell = "E2"; Value = "1.77%" }
    @{ Cell = "D3"; Value = "31.87" }
    @{ Cell = "E3"; Value = "0.70%" }
    @{ Cell = "D4"; Value = "4.980" }
    @{ Cell = "E4"; Value = "-2.02%" }
    @{ Cell = "D5"; Value = "0.07823" }
    @{ Cell = "E5"; Value = "-2.56%" }
    @{ Cell = "D6"; Value = "2.151" }
    @{ Cell = "E6"; Value = "-16.86%" }
    @{ Cell = "D7"; Value = "7.801" }
    @{ Cell = "E7"; Value = "0.04%" }
    @{ Cell = "D8"; Value = "3.786" }
    @{ Cell = "E8"; Value = "-0.88%" }
    @{ Cell = "D9"; Value = "0.9198" }
    @{ Cell = "E9"; Value = "-0.55%" }
    @{ Cell = "D10"; Value = "0.1754" }
    @{ Cell = "E10"; Value = "-0.37%" }
    @{ Cell = "D11"; Value = "0.07764" }
    @{ Cell = "E11"; Value = "4.94%" }
    @{ Cell = "D12"; Value = "0.08914" }
    @{ Cell = "E12"; Value = "-0.91%" }
    @{ Cell = "D13"; Value = "0.03086" }
    @{ Cell = "E13"; Value = "0.98%" }
    @{ Cell = "E14"; Value = "0.14%" }
    @{ Cell = "D15"; Value = "0.001506" }
    @{ Cell = "E15"; Value = "0.12%" }
    @{ Cell = "D16"; Value = "0.005917" }
    @{ Cell = "E16"; Value = "-2.01%" }
    @{ Cell = "D17"; Value = "3.461" }
    @{ Cell = "E17"; Value = "-2.17%" }
    @{ Cell = "D18"; Value = "2.267" }
    @{ Cell = "E18"; Value = "0.85%" }
    @{ Cell = "D20"; Value = "0.1328" }
    @{ Cell = "E20"; Value = "-0.60%" }
    @{ Cell = "D21"; Value = "4.157" }
    @{ Cell = "E21"; Value = "3.29%" }
    @{ Cell = "D22"; Value = "0.1793" }
    @{ Cell = "E22"; Value = "8.83%" }
    @{ Cell = "D23"; Value = "0.04590" }
    @{ Cell = "E23"; Value = "0.03%" }
    @{ Cell = "D24"; Value = "0.001241" }
    @{ Cell = "E24"; Value = "-0.21%" }
    @{ Cell = "D25"; Value = "0.004479" }
    @{ Cell = "E25"; Value = "1.14%" }
    @{ Cell = "D26"; Value = "0.0001250" }
    @{ Cell = "E26"; Value = "4.16%" }
    @{ Cell = "D39"; Value = "0.01769" }
    @{ Cell = "E39"; Value = "0.51%" }
    @{ Cell = "D40"; Value = "0.04771" }
    @{ Cell = "E40"; Value = "5.85%" }
    @{ Cell = "D41"; Value = "0.007121" }
    @{ Cell = "E41"; Value = "4.00%" }
    @{ Cell = "D42"; Value = "0.1370" }
    @{ Cell = "E42"; Value = "1.71%" }
    @{ Cell = "D43"; Value = "0.002150" }
    @{ Cell = "E43"; Value = "-2.72%" }
    @{ Cell = "D44"; Value = "0.01080" }
    @{ Cell = "E44"; Value = "9.76%" }
    @{ Cell = "D45"; Value = "0.00006235" }
    @{ Cell = "E45"; Value = "-3.47%" }
    @{ Cell = "E46"; Value = "0.00%" }
    @{ Cell = "E47"; Value = "-59.31%" }
    @{ Cell = "D48"; Value = "1.175" }
    @{ Cell = "E48"; Value = "43.15%" }
    @{ Cell = "E49"; Value = "0.00%" }
    @{ Cell = "E50"; Value = "0.00%" }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
}
